$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.113709449768066
$ws.Range("B1").Value = 2.823112964630127
$ws.Range("C1").Value = 1.956591606140137
$ws.Range("D1").Value = 1.751981139183044
$ws.Range("E1").Value = 1.67714536190033
